# Auto-generated PowerShell COM-interop script to apply market-data refresh edits
# to Sheets/Aegis_Profits.xlsx (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 288.2857
$ws.Range("I52").Value = 203
$ws.Range("J52").Value = 800
$ws.Range("K52").Value = 609
$ws.Range("L52").Value = 2400
$ws.Range("M52").Value = -449
$ws.Range("N52").Value = -2720
$ws.Range("H62").Value = 1500
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -876
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 1500
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -4380
$ws.Range("N65").ClearContents()
$ws.Range("H80").Value = 40447.2
$ws.Range("I80").Value = 240
$ws.Range("J80").Value = 100758
$ws.Range("K80").Value = 720
$ws.Range("L80").Value = 302274
$ws.Range("M80").Value = 278
$ws.Range("N80").Value = -304270
$ws.Range("H83").Value = 40447.2
$ws.Range("I83").Value = 240
$ws.Range("J83").Value = 100758
$ws.Range("K83").Value = 2160
$ws.Range("L83").Value = 906822
$ws.Range("M83").Value = 2832
$ws.Range("N83").Value = -916806
$ws.Range("H98").Value = 3371.25
$ws.Range("I98").Value = 3371.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 3371.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1873.25
$ws.Range("N98").ClearContents()
$ws.Range("H112").Value = 1184.8334
$ws.Range("J112").Value = 1184.8334
$ws.Range("L112").Value = 3554.5002
$ws.Range("N112").Value = -5770.5002
$ws.Range("H122").Value = 3371.25
$ws.Range("I122").Value = 3371.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10113.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7663.75
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 1391.6216
$ws.Range("I137").Value = 1140.3667
$ws.Range("J137").Value = 2468.4285
$ws.Range("K137").Value = 3421.1001
$ws.Range("L137").Value = 7405.2855
$ws.Range("M137").Value = -871.1001000000001
$ws.Range("N137").Value = -12505.2855
$ws.Range("H141").Value = 2534.2104
$ws.Range("I141").Value = 1947.1428
$ws.Range("J141").Value = 4178
$ws.Range("K141").Value = 5841.428400000001
$ws.Range("L141").Value = 12534
$ws.Range("M141").Value = -661.4284000000007
$ws.Range("N141").Value = -22894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 84782.664
$ws.Range("I45").Value = 112090.22
$ws.Range("K45").Value = 112090.22
$ws.Range("M45").Value = -111713.22
$ws.Range("H74").Value = 1290.2903
$ws.Range("I74").Value = 637.88464
$ws.Range("J74").Value = 4682.8
$ws.Range("K74").Value = 637.88464
$ws.Range("L74").Value = 4682.8
$ws.Range("M74").Value = 236.11536
$ws.Range("N74").Value = -6430.8
$ws.Range("H77").Value = 1290.2903
$ws.Range("I77").Value = 637.88464
$ws.Range("J77").Value = 4682.8
$ws.Range("K77").Value = 3189.4232
$ws.Range("L77").Value = 23414
$ws.Range("M77").Value = 1178.5768
$ws.Range("N77").Value = -32150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 325.25
$ws.Range("I22").Value = 267
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 267
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -94
$ws.Range("N22").Value = -846
$ws.Range("H94").Value = 800.38464
$ws.Range("I94").Value = 650.8333
$ws.Range("J94").Value = 928.5714
$ws.Range("K94").Value = 650.8333
$ws.Range("L94").Value = 928.5714
$ws.Range("M94").Value = -199.8333
$ws.Range("N94").Value = -1830.5714
$ws.Range("H134").Value = 2268.2708
$ws.Range("I134").Value = 2253.848
$ws.Range("J134").Value = 2600
$ws.Range("K134").Value = 6761.544
$ws.Range("L134").Value = 7800
$ws.Range("M134").Value = -4226.544
$ws.Range("N134").Value = -12870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26746.207
$ws.Range("I31").Value = 958.5833
$ws.Range("J31").Value = 42615.51
$ws.Range("K31").Value = 958.5833
$ws.Range("L31").Value = 42615.51
$ws.Range("M31").Value = -663.5833
$ws.Range("N31").Value = -43205.51
$ws.Range("H32").Value = 24905
$ws.Range("I32").Value = 20010
$ws.Range("K32").Value = 20010
$ws.Range("M32").Value = -19694
$ws.Range("H34").Value = 26746.207
$ws.Range("I34").Value = 958.5833
$ws.Range("J34").Value = 42615.51
$ws.Range("K34").Value = 958.5833
$ws.Range("L34").Value = 42615.51
$ws.Range("M34").Value = -756.5833
$ws.Range("N34").Value = -43019.51
$ws.Range("H134").Value = 1049.1818
$ws.Range("I134").Value = 953.7059
$ws.Range("J134").Value = 1373.8
$ws.Range("K134").Value = 2861.1177
$ws.Range("L134").Value = 4121.4
$ws.Range("M134").Value = -326.1177000000002
$ws.Range("N134").Value = -9191.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1558.5714
$ws.Range("I5").Value = 771.25
$ws.Range("J5").Value = 2221.5789
$ws.Range("K5").Value = 2313.75
$ws.Range("L5").Value = 6664.736699999999
$ws.Range("M5").Value = -2201.75
$ws.Range("N5").Value = -6888.736699999999
$ws.Range("H44").Value = 371.46667
$ws.Range("I44").Value = 248.625
$ws.Range("J44").Value = 511.85715
$ws.Range("K44").Value = 745.875
$ws.Range("L44").Value = 1535.57145
$ws.Range("M44").Value = -347.875
$ws.Range("N44").Value = -2331.57145
$ws.Range("H131").Value = 834.36
$ws.Range("I131").Value = 542.5
$ws.Range("J131").Value = 859.73914
$ws.Range("K131").Value = 1627.5
$ws.Range("L131").Value = 2579.21742
$ws.Range("M131").Value = 3412.5
$ws.Range("N131").Value = -12659.21742
$ws.Range("H132").Value = 2131.9443
$ws.Range("J132").Value = 2198.5293
$ws.Range("L132").Value = 19786.7637
$ws.Range("N132").Value = -24846.7637
$ws.Range("H135").Value = 1558.5714
$ws.Range("I135").Value = 771.25
$ws.Range("J135").Value = 2221.5789
$ws.Range("K135").Value = 6941.25
$ws.Range("L135").Value = 19994.2101
$ws.Range("M135").Value = -4406.25
$ws.Range("N135").Value = -25064.2101
$ws.Range("H137").Value = 18523520
$ws.Range("I137").Value = 9156
$ws.Range("J137").Value = 25644428
$ws.Range("K137").Value = 27468
$ws.Range("L137").Value = 76933284
$ws.Range("M137").Value = -22368
$ws.Range("N137").Value = -76943484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 147050.36
$ws.Range("I70").Value = 225632.89
$ws.Range("J70").Value = 5601.8
$ws.Range("K70").Value = 225632.89
$ws.Range("L70").Value = 5601.8
$ws.Range("M70").Value = -225362.89
$ws.Range("N70").Value = -6141.8
$ws.Range("H73").Value = 147050.36
$ws.Range("I73").Value = 225632.89
$ws.Range("J73").Value = 5601.8
$ws.Range("K73").Value = 225632.89
$ws.Range("L73").Value = 5601.8
$ws.Range("M73").Value = -224696.89
$ws.Range("N73").Value = -7473.8
$ws.Range("H126").Value = 3384.2
$ws.Range("I126").Value = 3204.6667
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9614.000100000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -7144.000100000001
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 1631.129
$ws.Range("I132").Value = 1233.2693
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 3699.8079
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -1169.8079
$ws.Range("N132").Value = -16160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 859.15
$ws.Range("J22").Value = 786.5
$ws.Range("L22").Value = 786.5
$ws.Range("N22").Value = -1376.5
$ws.Range("H27").Value = 859.15
$ws.Range("J27").Value = 786.5
$ws.Range("L27").Value = 786.5
$ws.Range("N27").Value = -1000.5
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 1666.6666
$ws.Range("K40").Value = 1666.6666
$ws.Range("M40").Value = -1530.6666
$ws.Range("H136").Value = 2053.0667
$ws.Range("I136").Value = 1984.3077
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 5952.9231
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -3402.9231
$ws.Range("N136").Value = -12600
